# Generate Report for handoff
# Updates the "Latest Handoff Datetime" column (D) for the
# "9c93ff0a-287c-4b6f-b9be-0f63c412a9e6.md" row (row 4) on both the
# zh-cn and de-de localization-status sheets, recording the timestamp
# of the new handoff that was just generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(4, 4).Value = "2016-01-15 09:57:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(4, 4).Value = "2016-01-15 09:58:08"
